$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect("D382")
$ws.Range("D2").Value = 0.4903322643479985
$ws.Protect("D382")
